$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2664359861591695
$ws.Range("C2").Value = 0.4429065743944637
$ws.Range("J2").Value = 0.02076124567474048
$ws.Range("P2").Value = 0.1799307958477509
$ws.Range("S2").Value = 0.08996539792387544
$ws.Range("C3").Value = 0.03424657534246575
$ws.Range("J3").Value = 0.0410958904109589
$ws.Range("P3").Value = 0.6917808219178082
$ws.Range("S3").Value = 0.2328767123287671
$ws.Range("J4").Value = 0.1388888888888889
$ws.Range("P4").Value = 0.6388888888888888
$ws.Range("S4").Value = 0.2222222222222222
$ws.Range("B6").Value = 0.07878787878787878
$ws.Range("D6").Value = 0.01818181818181818
$ws.Range("F6").Value = 0.01212121212121212
$ws.Range("J6").Value = 0.2727272727272727
$ws.Range("O6").Value = 0.01818181818181818
$ws.Range("Q6").Value = 0.2242424242424242
$ws.Range("R6").Value = 0.0303030303030303
$ws.Range("S6").Value = 0.3454545454545455
$ws.Range("B7").Value = 0.1883116883116883
$ws.Range("D7").Value = 0.01948051948051948
$ws.Range("E7").Value = 0.006493506493506494
$ws.Range("F7").Value = 0.05194805194805195
$ws.Range("J7").Value = 0.1363636363636364
$ws.Range("O7").Value = 0.01948051948051948
$ws.Range("Q7").Value = 0.1883116883116883
$ws.Range("R7").Value = 0.03246753246753246
$ws.Range("S7").Value = 0.3571428571428572
$ws.Range("B8").Value = 0.08226221079691516
$ws.Range("D8").Value = 0.012853470437018
$ws.Range("F8").Value = 0.07197943444730077
$ws.Range("J8").Value = 0.1182519280205656
$ws.Range("O8").Value = 0.02313624678663239
$ws.Range("Q8").Value = 0.1902313624678663
$ws.Range("R8").Value = 0.06683804627249357
$ws.Range("S8").Value = 0.4344473007712082
$ws.Range("B9").Value = 0.1061946902654867
$ws.Range("D9").Value = 0.008849557522123894
$ws.Range("F9").Value = 0.06194690265486726
$ws.Range("J9").Value = 0.1327433628318584
$ws.Range("O9").Value = 0.03539823008849557
$ws.Range("Q9").Value = 0.1592920353982301
$ws.Range("R9").Value = 0.07964601769911504
$ws.Range("S9").Value = 0.415929203539823
$ws.Range("B10").Value = 0.1102430555555556
$ws.Range("D10").Value = 0.01822916666666667
$ws.Range("F10").Value = 0.05295138888888889
$ws.Range("J10").Value = 0.1397569444444444
$ws.Range("O10").Value = 0.01909722222222222
$ws.Range("Q10").Value = 0.2204861111111111
$ws.Range("R10").Value = 0.08506944444444445
$ws.Range("S10").Value = 0.3541666666666667
$ws.Range("G11").Value = 0.1853281853281853
$ws.Range("J11").Value = 0.08494208494208494
$ws.Range("K11").Value = 0.2355212355212355
$ws.Range("L11").Value = 0.4787644787644788
$ws.Range("S11").Value = 0.01544401544401544
$ws.Range("G12").Value = 0.6533333333333333
$ws.Range("J12").Value = 0.2266666666666667
$ws.Range("K12").Value = 0.006666666666666667
$ws.Range("L12").Value = 0.04
$ws.Range("S12").Value = 0.07333333333333333
$ws.Range("G13").Value = 0.5806451612903226
$ws.Range("J13").Value = 0.3225806451612903
$ws.Range("S13").Value = 0.09677419354838709
$ws.Range("F15").Value = 0.0124223602484472
$ws.Range("H15").Value = 0.1366459627329193
$ws.Range("I15").Value = 0.08074534161490683
$ws.Range("J15").Value = 0.3975155279503105
$ws.Range("K15").Value = 0.03726708074534162
$ws.Range("M15").Value = 0.0124223602484472
$ws.Range("O15").Value = 0.04347826086956522
$ws.Range("S15").Value = 0.2795031055900621
$ws.Range("F16").Value = 0.01851851851851852
$ws.Range("H16").Value = 0.1790123456790123
$ws.Range("I16").Value = 0.08641975308641975
$ws.Range("J16").Value = 0.4259259259259259
$ws.Range("K16").Value = 0.1234567901234568
$ws.Range("M16").Value = 0.01851851851851852
$ws.Range("O16").Value = 0.0308641975308642
$ws.Range("S16").Value = 0.1172839506172839
$ws.Range("F17").Value = 0.01354401805869074
$ws.Range("H17").Value = 0.1196388261851016
$ws.Range("I17").Value = 0.1196388261851016
$ws.Range("J17").Value = 0.4469525959367946
$ws.Range("K17").Value = 0.1106094808126411
$ws.Range("M17").Value = 0.006772009029345372
$ws.Range("N17").Value = 0.002257336343115124
$ws.Range("O17").Value = 0.03160270880361174
$ws.Range("S17").Value = 0.1489841986455982
$ws.Range("F18").Value = 0.0131578947368421
$ws.Range("H18").Value = 0.1578947368421053
$ws.Range("I18").Value = 0.1776315789473684
$ws.Range("J18").Value = 0.4276315789473684
$ws.Range("K18").Value = 0.05921052631578947
$ws.Range("M18").Value = 0.006578947368421052
$ws.Range("N18").Value = 0.006578947368421052
$ws.Range("O18").Value = 0.0131578947368421
$ws.Range("S18").Value = 0.1381578947368421
$ws.Range("F19").Value = 0.01266891891891892
$ws.Range("H19").Value = 0.2119932432432433
$ws.Range("I19").Value = 0.09966216216216216
$ws.Range("J19").Value = 0.3378378378378378
$ws.Range("K19").Value = 0.09881756756756757
$ws.Range("M19").Value = 0.01858108108108108
$ws.Range("N19").Value = 0.001689189189189189
$ws.Range("O19").Value = 0.05658783783783784
$ws.Range("S19").Value = 0.1621621621621622
